$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Scenario:" table cell - append " - Sensor" right after "Blutzuckermessen"
#    "Blutzuckermessen " -> "Blutzuckermessen - Sensor "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Blutzuckermessen ", $true, $false, $false, $false, $false, $true, 1, $false, "Blutzuckermessen - Sensor ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Row 1.0 description: insert "in bestimmten Intervall " before
#    "Blutzuckergehalt im Blut"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" Blutzuckergehalt im Blut", $true, $false, $false, $false, $false, $true, 1, $false, " in bestimmten Intervall Blutzuckergehalt im Blut", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Row 1.1 description: replace "Insulin spritze" with "PIP"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Insulin spritze", $true, $false, $false, $false, $false, $true, 1, $false, "PIP", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Fill in the three previously-empty rows (1.4 / 1.5 / 1.6) of the second
#    table (Nr. / Actor / Description).
# ---------------------------------------------------------------------------
$t2 = $d.Tables.Item(2)

$row14 = $t2.Rows.Item(6)
$row14.Cells.Item(1).Range.Text = "1.4"
$row14.Cells.Item(2).Range.Text = "Sensor"
$row14.Cells.Item(3).Range.Text = "Bei Extremwerte Aufforderung zur manuellen Zweitmessung"

$row15 = $t2.Rows.Item(7)
$row15.Cells.Item(1).Range.Text = "1.5"
$row15.Cells.Item(2).Range.Text = "Sensor"
$row15.Cells.Item(3).Range.Text = "Nach Blutzuckermessung automatische Angabe der Insulinmenge, welche verabreicht werden muss"

$row16 = $t2.Rows.Item(8)
$row16.Cells.Item(1).Range.Text = "1.6"
$row16.Cells.Item(2).Range.Text = "Sensor"
$row16.Cells.Item(3).Range.Text = "Sensor sendet Messergebnisse an Gerät "

# ---------------------------------------------------------------------------
# 5) Remove the stray "_GoBack" bookmark in the third table (row 1.0.4).
#    The bookmark straddles "anderen " / "Blutzuckermessgerät " - a
#    find/replace that spans across that point rebuilds the run list and
#    drops the now-orphaned bookmark marks.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("anderen Blutzuckermessgerät", $true, $false, $false, $false, $false, $true, 1, $false, "anderen Blutzuckermessgerät", 2) | Out-Null
